$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.4444444444444444
$ws.Cells.Item(2, 3).Value = 0.3076923076923077
$ws.Cells.Item(2, 4).Value = 0.3636363636363637
$ws.Cells.Item(2, 5).Value = 13

$ws.Cells.Item(3, 2).Value = 0.4
$ws.Cells.Item(3, 3).Value = 0.5454545454545454
$ws.Cells.Item(3, 4).Value = 0.4615384615384615
$ws.Cells.Item(3, 5).Value = 11

$ws.Cells.Item(4, 2).Value = 0.4166666666666667
$ws.Cells.Item(4, 3).Value = 0.4166666666666667
$ws.Cells.Item(4, 4).Value = 0.4166666666666667
$ws.Cells.Item(4, 5).Value = 0.4166666666666667

$ws.Cells.Item(5, 2).Value = 0.4222222222222222
$ws.Cells.Item(5, 3).Value = 0.4265734265734266
$ws.Cells.Item(5, 4).Value = 0.4125874125874126
$ws.Cells.Item(5, 5).Value = 24

$ws.Cells.Item(6, 2).Value = 0.424074074074074
$ws.Cells.Item(6, 3).Value = 0.4166666666666667
$ws.Cells.Item(6, 4).Value = 0.4085081585081585
$ws.Cells.Item(6, 5).Value = 24

$ws.Cells.Item(7, 2).Value = 0.3
$ws.Cells.Item(7, 3).Value = 0.2307692307692308
$ws.Cells.Item(7, 4).Value = 0.2608695652173913
$ws.Cells.Item(7, 5).Value = 13

$ws.Cells.Item(8, 2).Value = 0.2857142857142857
$ws.Cells.Item(8, 3).Value = 0.3636363636363636
$ws.Cells.Item(8, 4).Value = 0.32
$ws.Cells.Item(8, 5).Value = 11

$ws.Cells.Item(9, 2).Value = 0.2916666666666667
$ws.Cells.Item(9, 3).Value = 0.2916666666666667
$ws.Cells.Item(9, 4).Value = 0.2916666666666667
$ws.Cells.Item(9, 5).Value = 0.2916666666666667

$ws.Cells.Item(10, 2).Value = 0.2928571428571428
$ws.Cells.Item(10, 3).Value = 0.2972027972027972
$ws.Cells.Item(10, 4).Value = 0.2904347826086957
$ws.Cells.Item(10, 5).Value = 24

$ws.Cells.Item(11, 2).Value = 0.293452380952381
$ws.Cells.Item(11, 3).Value = 0.2916666666666667
$ws.Cells.Item(11, 4).Value = 0.2879710144927536
$ws.Cells.Item(11, 5).Value = 24

$ws.Cells.Item(12, 2).Value = 0.5882352941176471
$ws.Cells.Item(12, 3).Value = 0.7692307692307693
$ws.Cells.Item(12, 4).Value = 0.6666666666666667
$ws.Cells.Item(12, 5).Value = 13

$ws.Cells.Item(13, 2).Value = 0.5714285714285714
$ws.Cells.Item(13, 3).Value = 0.3636363636363636
$ws.Cells.Item(13, 4).Value = 0.4444444444444444
$ws.Cells.Item(13, 5).Value = 11

$ws.Cells.Item(14, 2).Value = 0.5833333333333334
$ws.Cells.Item(14, 3).Value = 0.5833333333333334
$ws.Cells.Item(14, 4).Value = 0.5833333333333334
$ws.Cells.Item(14, 5).Value = 0.5833333333333334

$ws.Cells.Item(15, 2).Value = 0.5798319327731092
$ws.Cells.Item(15, 3).Value = 0.5664335664335665
$ws.Cells.Item(15, 4).Value = 0.5555555555555556
$ws.Cells.Item(15, 5).Value = 24

$ws.Cells.Item(16, 2).Value = 0.5805322128851541
$ws.Cells.Item(16, 3).Value = 0.5833333333333334
$ws.Cells.Item(16, 4).Value = 0.5648148148148149
$ws.Cells.Item(16, 5).Value = 24

$ws.Cells.Item(17, 2).Value = 0.4166666666666667
$ws.Cells.Item(17, 3).Value = 0.3846153846153846
$ws.Cells.Item(17, 4).Value = 0.4
$ws.Cells.Item(17, 5).Value = 13

$ws.Cells.Item(18, 2).Value = 0.3333333333333333
$ws.Cells.Item(18, 3).Value = 0.3636363636363636
$ws.Cells.Item(18, 4).Value = 0.3478260869565217
$ws.Cells.Item(18, 5).Value = 11

$ws.Cells.Item(19, 2).Value = 0.375
$ws.Cells.Item(19, 3).Value = 0.375
$ws.Cells.Item(19, 4).Value = 0.375
$ws.Cells.Item(19, 5).Value = 0.375

$ws.Cells.Item(20, 2).Value = 0.375
$ws.Cells.Item(20, 3).Value = 0.3741258741258742
$ws.Cells.Item(20, 4).Value = 0.3739130434782609
$ws.Cells.Item(20, 5).Value = 24

$ws.Cells.Item(21, 2).Value = 0.3784722222222223
$ws.Cells.Item(21, 3).Value = 0.375
$ws.Cells.Item(21, 4).Value = 0.3760869565217391
$ws.Cells.Item(21, 5).Value = 24

$ws.Cells.Item(22, 2).Value = 0.6363636363636364
$ws.Cells.Item(22, 3).Value = 0.5384615384615384
$ws.Cells.Item(22, 4).Value = 0.5833333333333334
$ws.Cells.Item(22, 5).Value = 13

$ws.Cells.Item(23, 2).Value = 0.5384615384615384
$ws.Cells.Item(23, 3).Value = 0.6363636363636364
$ws.Cells.Item(23, 4).Value = 0.5833333333333334
$ws.Cells.Item(23, 5).Value = 11

$ws.Cells.Item(24, 2).Value = 0.5833333333333334
$ws.Cells.Item(24, 3).Value = 0.5833333333333334
$ws.Cells.Item(24, 4).Value = 0.5833333333333334
$ws.Cells.Item(24, 5).Value = 0.5833333333333334

$ws.Cells.Item(25, 2).Value = 0.5874125874125874
$ws.Cells.Item(25, 3).Value = 0.5874125874125874
$ws.Cells.Item(25, 4).Value = 0.5833333333333334
$ws.Cells.Item(25, 5).Value = 24

$ws.Cells.Item(26, 2).Value = 0.5914918414918415
$ws.Cells.Item(26, 3).Value = 0.5833333333333334
$ws.Cells.Item(26, 4).Value = 0.5833333333333334
$ws.Cells.Item(26, 5).Value = 24
